$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells: D1 = "x-25", E1 = "y-25"
$ws.Range("D1").Value = "x-25"
$ws.Range("E1").Value = "y-25"

# New formula columns: D = B-25, E = C-25, for every data row (2..62)
for ($r = 2; $r -le 62; $r++) {
    $ws.Cells.Item($r, 4).Formula = "=B$r-25"
    $ws.Cells.Item($r, 5).Formula = "=C$r-25"
}

# Highlight column D (yellow) and column E (blue), header + all data rows
$ws.Range("D1:D62").Interior.Color = 65535
$ws.Range("E1:E62").Interior.Color = 15773696

# Match the author's final selection
$ws.Range("E15").Select()
